# Update Saldo_guide worksheet:
#  - Column G ("Dt. Referencia") moves forward one day (45407 -> 45408) for every data row
#  - A handful of rows get revised "Saldo Previsto" (D), "Vl. Projetado" (E) and
#    "Vl. Total" (H) figures

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saldo_guide")

# Every data row (2-310) has its reference date bumped one day forward
# (45407 = 2024-04-25  ->  45408 = 2024-04-26).
$ws.Range("G2:G310").Value = 45408

# Rows whose Saldo Previsto / Vl. Projetado / Vl. Total values changed.
$updates = @(
    @{ Row = 6;   D = 0;        E = 0;        H = 0 },
    @{ Row = 12;  D = 26589.35; E = -10615.62; H = 15973.73 },
    @{ Row = 71;  D = 0;        E = 0;        H = 0 },
    @{ Row = 113; D = 0;        E = 0;        H = 0 },
    @{ Row = 124; D = 55551.16; E = 0;        H = 55551.16 },
    @{ Row = 151; D = 0;        E = 0;        H = 0 },
    @{ Row = 184; D = 26295.8;  E = -1687.24; H = 24608.56 },
    @{ Row = 290; D = 785.7;    E = 0;        H = 785.7 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}
